# Rename feature/response headers to short code-friendly names and tidy up
# sheet views, in preparation for using this workbook with the
# auto_model_reduction example notebook.

$wb = $excel.ActiveWorkbook

$trainData = $wb.Worksheets.Item("Train Data")
$testData  = $wb.Worksheets.Item("Test Data")
$design    = $wb.Worksheets.Item("Design Parameters")
$responses = $wb.Worksheets.Item("Responses")
$misc      = $wb.Worksheets.Item("Misc")

# --- Train Data -------------------------------------------------------------
$trainData.Range("B1").Value = "roll_temperature"
$trainData.Range("C1").Value = "density"
$trainData.Range("D1").Value = "porosity"
$trainData.Range("E1").Value = "C_5_6"
$trainData.Range("F1").Value = "_10C"
$trainData.Rows.Item(1).RowHeight = 30

# --- Test Data ----------------------------------------------------------------
$testData.Range("B1").Value = "roll_temperature"
$testData.Range("C1").Value = "density"
$testData.Range("D1").Value = "porosity"
$testData.Range("E1").Value = "C_5_6"
$testData.Range("F1").Value = "_10C"
$testData.Rows.Item(1).RowHeight = 30

# --- Design Parameters --------------------------------------------------------
$design.Range("B2").Value = "roll_temperature"
$design.Range("B3").Value = "density"
$design.Range("B4").Value = "porosity"

# B2:B4 should pick up the wrap/centered-without-vertical-center style that
# is already used for the header row of Train/Test Data (B1 there).
$trainData.Range("B1").Copy()
$design.Range("B2:B4").PasteSpecial(-4122)

# New (currently blank) cells further along row 6, formatted to match
# existing styles elsewhere in the workbook.
$trainData.Range("B1").Copy()
$design.Range("I6:K6").PasteSpecial(-4122)
$trainData.Range("A1").Copy()
$design.Range("L6:M6").PasteSpecial(-4122)

# --- Responses -----------------------------------------------------------------
$responses.Range("A2").Value = "C_5_6"
$responses.Range("A3").Value = "_10C"

# --- Misc ------------------------------------------------------------------
$misc.Range("A1").Value = "model"
$misc.Range("B1").Value = "Linear"

# --- Selections / active sheet ----------------------------------------------
$trainData.Activate()
$trainData.Range("F1").Select() | Out-Null

$testData.Activate()
$testData.Range("F1").Select() | Out-Null

$design.Activate()
$design.Range("O9").Select() | Out-Null

$misc.Activate()
$misc.Range("D4:J33").Select() | Out-Null
$misc.Cells.Item(4, 10).Activate() | Out-Null

$responses.Activate()
$responses.Range("A3").Select() | Out-Null
